$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 7088
$ws.Range("B2").Value = 63
$ws.Range("C2").Value = 2208
$ws.Range("D2").Value = 18
